$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-16 Saturday" "2025-08-17 Sunday"

Replace-Text "248×2=496" "599×4=2396"
Replace-Text "644×6=3864" "562×6=3372"
Replace-Text "876×2=1752" "506×8=4048"
Replace-Text "637×5=3185" "588×9=5292"
Replace-Text "883×4=3532" "135×3=405"

Replace-Text "719×5=3595" "268×4=1072"
Replace-Text "451×4=1804" "292×6=1752"
Replace-Text "340×9=3060" "933×6=5598"
Replace-Text "308×6=1848" "240×9=2160"
Replace-Text "123×3=369" "901×2=1802"

Replace-Text "954×8=7632" "593×6=3558"
Replace-Text "929×7=6503" "510×6=3060"
Replace-Text "440×8=3520" "935×4=3740"
Replace-Text "568×3=1704" "684×2=1368"
Replace-Text "602×4=2408" "907×2=1814"

Replace-Text "554×3=1662" "853×9=7677"
Replace-Text "308×7=2156" "514×9=4626"
Replace-Text "248×8=1984" "545×9=4905"
Replace-Text "360×2=720" "243×4=972"
Replace-Text "370×5=1850" "941×7=6587"

Replace-Text "919×9=8271" "745×3=2235"
Replace-Text "271×8=2168" "444×4=1776"
Replace-Text "756×3=2268" "606×2=1212"
Replace-Text "771×6=4626" "829×6=4974"
Replace-Text "204×7=1428" "801×6=4806"
